$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - update MCC, Accuracy, Balanced Accuracy, F1
$ws.Range("E2").Value = 0.96064535921058791
$ws.Range("F2").Value = 0.98
$ws.Range("G2").Value = 0.9814814814814814
$ws.Range("H2").Value = 0.98113207547169812

# Row 10 - update Threshold, MCC, Accuracy, Balanced Accuracy, F1
$ws.Range("D10").Value = 0.9
$ws.Range("E10").Value = 0.80336173105394182
$ws.Range("F10").Value = 0.97
$ws.Range("G10").Value = 0.83333333333333326
$ws.Range("H10").Value = 0.8

# Row 11 - update Threshold, MCC, Balanced Accuracy, F1
$ws.Range("D11").Value = 0.96
$ws.Range("E11").Value = 0.93767107939836591
$ws.Range("G11").Value = 0.99456521739130432
$ws.Range("H11").Value = 0.94117647058823528

# Clear the cell selection left over from previous editing session
$ws.Range("A1").Select()
